# Add "2022-Q3" data:
#  1. Insert a new worksheet "2022-Q3" right after "总计", before "2022-Q2".
#  2. Fill it with the fund-holding breakdown for 2022-Q3.
#  3. Insert a new summary row in "总计" for the 2022-Q3 totals, pushing the
#     existing rows down by one.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new "2022-Q3" sheet right after "总计" --------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Reuse the same header / column formatting as the existing quarter sheets
# (bold + centered + bordered header row, bold centered index column).
$q2Sheet.Range("A1:H7").Copy() | Out-Null
$newSheet.Range("A1").PasteSpecial(-4122) | Out-Null

# --- 2. Populate header row -------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- 3. Populate data rows ---------------------------------------------------
# Columns B, C, D, E, F, G are stored as text (even the numeric-looking ones),
# only A (index) and H (rank) are real numbers.
function Set-TextCell($range) {
    $range.NumberFormat = "@"
}

$data = @(
    @{ A=0; B="002121"; C="广发沪港深新起点股票A";                         D="25.97"; E="85.60"; F="4.90"; G="1.2725"; H=5 },
    @{ A=1; B="009896"; C="广发港股通成长精选股票A";                       D="18.30"; E="90.12"; F="5.18"; G="0.9479"; H=6 },
    @{ A=2; B="009897"; C="广发港股通成长精选股票C";                       D="5.86";  E="90.12"; F="5.18"; G="0.3035"; H=6 },
    @{ A=3; B="501021"; C="华宝标普香港上市中国中小盘指数（LOF）A";        D="4.19";  E="92.99"; F="2.08"; G="0.0872"; H=4 },
    @{ A=4; B="010024"; C="广发沪港深新起点股票C";                         D="0.49";  E="85.60"; F="4.90"; G="0.0240"; H=5 },
    @{ A=5; B="006127"; C="华宝标普香港上市中国中小盘指数（LOF）C";        D="0.24";  E="92.99"; F="2.08"; G="0.0050"; H=4 }
)

$r = 2
foreach ($row in $data) {
    $rangeB = $newSheet.Range("B$r`:G$r")
    Set-TextCell $rangeB

    $newSheet.Range("A$r").Value = $row.A
    $newSheet.Range("B$r").Value = $row.B
    $newSheet.Range("C$r").Value = $row.C
    $newSheet.Range("D$r").Value = $row.D
    $newSheet.Range("E$r").Value = $row.E
    $newSheet.Range("F$r").Value = $row.F
    $newSheet.Range("G$r").Value = $row.G
    $newSheet.Range("H$r").Value = $row.H

    # Drop the "Text" number-format style added above so the cell keeps the
    # plain/default style (only the type stays text), matching the rest of
    # the workbook's body cells.
    $newSheet.Range("B$r`:G$r").Style = "Normal"

    $r = $r + 1
}

# --- 4. Insert the new 2022-Q3 row into the "总计" summary sheet -----------
$summary.Rows.Item(2).Insert() | Out-Null

# Recreate the bold/centered/bordered index-column style for the new A2 cell
# by copying formats from A3 (identical style used by every row in column A).
$summary.Range("A3").Copy() | Out-Null
$summary.Range("A2").PasteSpecial(-4122) | Out-Null
$summary.Range("B2:D2").ClearFormats() | Out-Null

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 2.64
